$wb = $excel.ActiveWorkbook

# Sheets that contain the Pattern3-Data+News / deepseek-v3 performance row:
#  - "Summary" sheet, row 12
#  - "Pattern3-Data+News" sheet, row 2
$targets = @(
    @{ Sheet = "Summary"; Row = 12 },
    @{ Sheet = "Pattern3-Data+News"; Row = 2 }
)

foreach ($t in $targets) {
    $ws = $wb.Worksheets.Item($t.Sheet)
    $r = $t.Row

    # Columns D,E,F,G,J,K,L,P hold formatted numbers/dates stored as plain
    # text (inlineStr) in the workbook. Force text format while assigning so
    # Excel does not auto-convert the strings back into numeric/percentage
    # values, then restore the "Normal" style so no stray cell formatting
    # is left behind.
    $textVals = [ordered]@{
        "D" = "¥1,005,052.00"
        "E" = "¥+5,052.00"
        "F" = "+0.51%"
        "G" = "+23.36%"
        "J" = "60.0%"
        "K" = "0.1009%"
        "L" = "0.0840%"
        "P" = "20251224"
    }

    foreach ($col in $textVals.Keys) {
        $cell = $ws.Range("$col$r")
        $cell.NumberFormat = "@"
        $cell.Value = $textVals[$col]
        $cell.Style = "Normal"
    }

    $ws.Range("H$r").Value = 17.501
    $ws.Range("M$r").Value = 6
    $ws.Range("N$r").Value = 6
}
